$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 603
$ws1.Range("F8").Value = 130
$ws1.Range("F9").Value = 9058
$ws1.Range("F11").Value = 334
$ws1.Range("F12").Value = 1162
$ws1.Range("F13").Value = 1037
$ws1.Range("F14").Value = 128
$ws1.Range("F18").Value = 327
$ws1.Range("F20").Value = 240
$ws1.Range("F21").Value = 1162

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 603
$ws4.Range("F10").Value = 130
$ws4.Range("F11").Value = 9058
$ws4.Range("F13").Value = 334
$ws4.Range("F14").Value = 1162
$ws4.Range("F15").Value = 1037
$ws4.Range("F16").Value = 128
$ws4.Range("F20").Value = 327
$ws4.Range("F22").Value = 240
$ws4.Range("F23").Value = 1162
